$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New GSW box-score rows (142-145 -> sheet rows 144-147)
$data = @(
    @(142,"GSW","MIA","away","2025-03-25","240:00",33,83,0.398,9,38,0.237,11,16,0.6879999999999999,14,24,38,20,8,1,12,14,86,-26),
    @(143,"MIA","GSW","home","2025-03-25","240:00",43,77,0.5580000000000001,17,25,0.68,9,12,0.75,8,34,42,26,8,5,15,15,112,26),
    @(144,"GSW","NOP","away","2025-03-28","240:00",38,92,0.413,13,55,0.236,22,29,0.759,18,32,50,29,10,1,10,22,111,16),
    @(145,"NOP","GSW","home","2025-03-28","240:00",33,81,0.407,7,23,0.304,22,28,0.786,14,32,46,24,6,5,15,23,95,-16)
)

$startRow = 144
$dateCol = 5   # column E holds the DATE text, keep it as text, not an auto-converted date

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        if (($c + 1) -eq $dateCol) {
            # force text storage so "2025-03-25" style strings aren't reinterpreted as dates
            $cell.NumberFormat = "@"
            $cell.Value = $rowVals[$c]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowVals[$c]
        }
    }
}

# Column A on these data rows uses the same bold/bordered/centered style as the rest of
# the table's row-index column; copy that formatting from the row right above (A143).
$fmtSrc = $ws.Range("A143")
$fmtSrc.Copy()
$fmtDst = $ws.Range("A144:A147")
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = 0
